# Updates the cryptos list (coinranking.com snapshot) with refreshed
# Price (column D) and Volume/1h change (column E) values for most rows.
# Rows 43-45 (Stacks / Monero / WEMIXToken) are also re-ordered: their
# Coin name, Link, Price and Volume(1h) values rotate between the three
# rows while the row numbers/ranks stay fixed.
#
# D-column price cells are plain numeric-looking text (e.g. "352.83")
# in the source workbook, so we force each target cell to Text format
# before assigning the value -- otherwise Excel's COM layer would
# auto-convert the literal string into a floating point number and the
# exact textual representation (e.g. trailing zeros) would be lost.
# The cell style is reset back to "Normal" right after so no stray
# number-format style is left attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value2 = '51.658.78'
$c.Style = "Normal"
$ws.Range("E2").Value2 = '  +0.29%  '
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value2 = '2.825.24'
$c.Style = "Normal"
$ws.Range("E3").Value2 = '  +2.30%  '
$ws.Range("E4").Value2 = '  +0.07%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value2 = '352.83'
$c.Style = "Normal"
$ws.Range("E5").Value2 = '  +5.98%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value2 = '113.18'
$c.Style = "Normal"
$ws.Range("E6").Value2 = '  -2.54%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value2 = '0.573'
$c.Style = "Normal"
$ws.Range("E7").Value2 = '  +6.45%  '
$ws.Range("E8").Value2 = '  +0.02%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value2 = '0.596'
$c.Style = "Normal"
$ws.Range("E9").Value2 = '  +3.83%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value2 = '41.44'
$c.Style = "Normal"
$ws.Range("E10").Value2 = '  -0.83%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value2 = '0.0853'
$c.Style = "Normal"
$ws.Range("E11").Value2 = '  -1.70%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value2 = '20.02'
$c.Style = "Normal"
$ws.Range("E12").Value2 = '  -1.45%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value2 = '7.70'
$c.Style = "Normal"
$ws.Range("E14").Value2 = '  +0.53%  '
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value2 = '3.278.26'
$c.Style = "Normal"
$ws.Range("E15").Value2 = '  +2.69%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value2 = '2.830.84'
$c.Style = "Normal"
$ws.Range("E16").Value2 = '  +2.11%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value2 = '0.893'
$c.Style = "Normal"
$ws.Range("E17").Value2 = '  +0.34%  '
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value2 = '51.486.94'
$c.Style = "Normal"
$ws.Range("E18").Value2 = '  -0.08%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value2 = '7.34'
$c.Style = "Normal"
$ws.Range("E19").Value2 = '  +7.11%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value2 = '3.15'
$c.Style = "Normal"
$ws.Range("E20").Value2 = '  -4.01%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value2 = '13.40'
$c.Style = "Normal"
$ws.Range("E21").Value2 = '  -0.78%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value2 = '0.0₃0993'
$c.Style = "Normal"
$ws.Range("E22").Value2 = '  +1.73%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value2 = '270.60'
$c.Style = "Normal"
$ws.Range("E23").Value2 = '  -2.70%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value2 = '69.62'
$c.Style = "Normal"
$ws.Range("E24").Value2 = '  -0.11%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value2 = '2.75'
$c.Style = "Normal"
$ws.Range("E25").Value2 = '  +2.47%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value2 = '26.64'
$c.Style = "Normal"
$ws.Range("E26").Value2 = '  -0.70%  '
$ws.Range("E27").Value2 = '  -0.03%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value2 = '10.29'
$c.Style = "Normal"
$ws.Range("E28").Value2 = '  +1.19%  '
$ws.Range("E29").Value2 = '  +1.01%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value2 = '0.138'
$c.Style = "Normal"
$ws.Range("E30").Value2 = '  -2.27%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value2 = '34.14'
$c.Style = "Normal"
$ws.Range("E31").Value2 = '  -2.62%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value2 = '50.68'
$c.Style = "Normal"
$ws.Range("E32").Value2 = '  +1.16%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value2 = '5.80'
$c.Style = "Normal"
$ws.Range("E33").Value2 = '  +4.25%  '
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value2 = '0.0441'
$c.Style = "Normal"
$ws.Range("E34").Value2 = '  +24.93%  '
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value2 = '0.0823'
$c.Style = "Normal"
$ws.Range("E35").Value2 = '  +0.20%  '
$ws.Range("E36").Value2 = '  -0.03%  '
$ws.Range("E37").Value2 = '  -0.81%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value2 = '3.18'
$c.Style = "Normal"
$ws.Range("E39").Value2 = '  -1.69%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value2 = '18.03'
$c.Style = "Normal"
$ws.Range("E40").Value2 = '  -4.88%  '
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value2 = '23.72'
$c.Style = "Normal"
$ws.Range("E41").Value2 = '  +2.56%  '
$ws.Range("E42").Value2 = '  +2.81%  '
$ws.Range("B43").Value2 = 'Monero'
$ws.Range("C43").Value2 = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value2 = '126.38'
$c.Style = "Normal"
$ws.Range("E43").Value2 = '  -0.80%  '
$ws.Range("B44").Value2 = 'WEMIXToken'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value2 = '2.30'
$c.Style = "Normal"
$ws.Range("E44").Value2 = '  -0.15%  '
$ws.Range("B45").Value2 = 'Stacks'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value2 = '2.50'
$c.Style = "Normal"
$ws.Range("E45").Value2 = '  +1.39%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value2 = '2.076.89'
$c.Style = "Normal"
$ws.Range("E46").Value2 = '  -0.50%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value2 = '3.31'
$c.Style = "Normal"
$ws.Range("E47").Value2 = '  +0.02%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value2 = '5.68'
$c.Style = "Normal"
$ws.Range("E49").Value2 = '  +2.70%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value2 = '0.917'
$c.Style = "Normal"
$ws.Range("E50").Value2 = '  +4.99%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value2 = '60.66'
$c.Style = "Normal"
$ws.Range("E51").Value2 = '  +0.92%  '
